$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 6830.4
$ws.Range("I12").Value = 7798
$ws.Range("J12").Value = 541
$ws.Range("K12").Value = 7798
$ws.Range("L12").Value = 541
$ws.Range("M12").Value = -7628
$ws.Range("N12").Value = -881
# Row 17
$ws.Range("H17").Value = 261.31708
$ws.Range("J17").Value = 265.35
$ws.Range("L17").Value = 796.0500000000001
$ws.Range("N17").Value = -1132.05
# Row 135
$ws.Range("H135").Value = 957.75
$ws.Range("I135").Value = 978.6667
$ws.Range("J135").Value = 895
$ws.Range("K135").Value = 8808.0003
$ws.Range("L135").Value = 8055
$ws.Range("M135").Value = -6273.0003
$ws.Range("N135").Value = -13125
# Row 137
$ws.Range("H137").Value = 6675.524
$ws.Range("I137").Value = 11448.4
$ws.Range("J137").Value = 2336.5454
$ws.Range("K137").Value = 34345.2
$ws.Range("L137").Value = 7009.6362
$ws.Range("M137").Value = -31795.2
$ws.Range("N137").Value = -12109.6362

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 11
$ws.Range("H11").Value = 5010001.5
$ws.Range("I11").Value = 5010001.5
$ws.Range("K11").Value = 5010001.5
$ws.Range("M11").Value = -5009857.5
# Row 32
$ws.Range("H32").Value = 10952.558
$ws.Range("I32").Value = 6697.282
$ws.Range("J32").Value = 23718.385
$ws.Range("K32").Value = 6697.282
$ws.Range("L32").Value = 23718.385
$ws.Range("M32").Value = -6410.282
$ws.Range("N32").Value = -24292.385
# Row 45
$ws.Range("H45").Value = 68391.8
$ws.Range("I45").Value = 92378.45
$ws.Range("J45").Value = 2428.5
$ws.Range("K45").Value = 92378.45
$ws.Range("L45").Value = 2428.5
$ws.Range("M45").Value = -92001.45
$ws.Range("N45").Value = -3182.5
# Row 61
$ws.Range("H61").Value = 7754571
$ws.Range("I61").Value = 12822724
$ws.Range("K61").Value = 12822724
$ws.Range("M61").Value = -12822512
# Row 74
$ws.Range("H74").Value = 13048969
$ws.Range("I74").Value = 17647698
$ws.Range("J74").Value = 19233.166
$ws.Range("K74").Value = 17647698
$ws.Range("L74").Value = 19233.166
$ws.Range("M74").Value = -17646824
$ws.Range("N74").Value = -20981.166
# Row 77
$ws.Range("H77").Value = 13048969
$ws.Range("I77").Value = 17647698
$ws.Range("J77").Value = 19233.166
$ws.Range("K77").Value = 88238490
$ws.Range("L77").Value = 96165.83
$ws.Range("M77").Value = -88234122
$ws.Range("N77").Value = -104901.83
# Row 136
$ws.Range("H136").Value = 7754571
$ws.Range("I136").Value = 12822724
$ws.Range("K136").Value = 38468172
$ws.Range("M136").Value = -38465622

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1931.3572
$ws.Range("I31").Value = 1853
$ws.Range("J31").Value = 1957.4762
$ws.Range("K31").Value = 1853
$ws.Range("L31").Value = 1957.4762
$ws.Range("M31").Value = -1558
$ws.Range("N31").Value = -2547.4762
# Row 34
$ws.Range("H34").Value = 1931.3572
$ws.Range("I34").Value = 1853
$ws.Range("J34").Value = 1957.4762
$ws.Range("K34").Value = 1853
$ws.Range("L34").Value = 1957.4762
$ws.Range("M34").Value = -1651
$ws.Range("N34").Value = -2361.4762
# Row 58
$ws.Range("H58").Value = 1867.3043
$ws.Range("I58").Value = 1714.6666
$ws.Range("J58").Value = 2153.5
$ws.Range("K58").Value = 1714.6666
$ws.Range("L58").Value = 2153.5
$ws.Range("M58").Value = -1511.6666
$ws.Range("N58").Value = -2559.5
# Row 94
$ws.Range("H94").Value = 983.6316
$ws.Range("I94").Value = 1204
$ws.Range("J94").Value = 942.3125
$ws.Range("K94").Value = 1204
$ws.Range("L94").Value = 942.3125
$ws.Range("M94").Value = -753
$ws.Range("N94").Value = -1844.3125
# Row 136
$ws.Range("H136").Value = 1867.3043
$ws.Range("I136").Value = 1714.6666
$ws.Range("J136").Value = 2153.5
$ws.Range("K136").Value = 5143.9998
$ws.Range("L136").Value = 6460.5
$ws.Range("M136").Value = -2593.9998
$ws.Range("N136").Value = -11560.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 110
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 3000
$ws.Range("M110").Value = 1090
# Row 132
$ws.Range("H132").Value = 111112380
$ws.Range("J132").Value = 938
$ws.Range("L132").Value = 8442
$ws.Range("N132").Value = -13502

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 3400
$ws.Range("I10").Value = 2600
$ws.Range("K10").Value = 2600
$ws.Range("M10").Value = -2460
# Row 132
$ws.Range("H132").Value = 4944.8887
$ws.Range("I132").Value = 5147.0264
$ws.Range("K132").Value = 15441.0792
$ws.Range("M132").Value = -12911.0792
# Row 136
$ws.Range("H136").Value = 3093.2354
$ws.Range("I136").Value = 2001.5
$ws.Range("J136").Value = 4063.6667
$ws.Range("K136").Value = 6004.5
$ws.Range("L136").Value = 12191.0001
$ws.Range("M136").Value = -3454.5
$ws.Range("N136").Value = -17291.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
# Row 58
$ws.Range("H58").Value = 14998
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
# Row 132
$ws.Range("H132").Value = 6088.892
$ws.Range("I132").Value = 8082.125
$ws.Range("J132").Value = 2409.077
$ws.Range("K132").Value = 24246.375
$ws.Range("L132").Value = 7227.231000000001
$ws.Range("M132").Value = -21716.375
$ws.Range("N132").Value = -12287.231
# Row 136
$ws.Range("H136").Value = 32000.846
$ws.Range("I136").Value = 10828.1
$ws.Range("J136").Value = 54287.95
$ws.Range("K136").Value = 32484.3
$ws.Range("L136").Value = 162863.85
$ws.Range("M136").Value = -29934.3
$ws.Range("N136").Value = -167963.85
